# This script re-applies a set of "findings" (rows) that had been
# accidentally shuffled between adjacent/nearby records. Each affected
# row's entire contents (columns A through AY) are exchanged with
# another row's contents so that the correct Id (column A) once again
# lines up with the correct record data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$FirstCol = 1    # A
$LastCol  = 51   # AY

function Get-CellValue($row, $col) {
    return $ws.Cells.Item($row, $col).Value2
}

function Set-CellValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    if ($null -eq $val) {
        $cell.Value2 = $null
    } elseif ($val -is [string]) {
        if ($val -eq "") {
            $cell.Value2 = ""
        } else {
            # Prefix with an apostrophe so that date/time/number-looking
            # text (e.g. "2026-01-13", "09:53") is kept as plain text
            # instead of being auto-converted by Excel into a date/number
            # serial value.
            $cell.Value2 = "'" + $val
        }
    } else {
        $cell.Value2 = $val
    }
}

function Get-RowValues($row) {
    $n = $LastCol - $FirstCol + 1
    $vals = New-Object 'object[]' $n
    for ($c = $FirstCol; $c -le $LastCol; $c++) {
        $vals[$c - $FirstCol] = Get-CellValue $row $c
    }
    return $vals
}

function Set-RowValues($row, $vals) {
    for ($c = $FirstCol; $c -le $LastCol; $c++) {
        Set-CellValue $row $c $vals[$c - $FirstCol]
    }
}

function Swap-Rows($r1, $r2) {
    $v1 = Get-RowValues $r1
    $v2 = Get-RowValues $r2
    Set-RowValues $r1 $v2
    Set-RowValues $r2 $v1
}

# Simple pairwise swaps: the two rows' data had been transposed.
Swap-Rows 29 30
Swap-Rows 36 37
Swap-Rows 48 49
Swap-Rows 50 51
Swap-Rows 60 61

# Rows 52, 53 and 54 form a 3-way rotation:
#   new row52 <- old row54
#   new row53 <- old row52
#   new row54 <- old row53
$v52 = Get-RowValues 52
$v53 = Get-RowValues 53
$v54 = Get-RowValues 54

Set-RowValues 52 $v54
Set-RowValues 53 $v52
Set-RowValues 54 $v53
